# Wimbledon ScoPo bracket update: record match winners for round 1 (rows
# 4-67) and round 2 (rows 70-133) by setting the "K" column indicator used
# by the sheet's VLOOKUP/IF formulas (1 = top player won, 2 = bottom player
# won). All downstream cells (L:T on ScoPo, and the "Results Export" sheet)
# are formula-driven and recalc automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ScoPo")

$winners = @{
    10 = 1
    22 = 2
    30 = 1
    31 = 2
    33 = 1
    37 = 2
    40 = 2
    41 = 2
    43 = 2
    44 = 1
    45 = 1
    46 = 1
    47 = 2
    48 = 2
    49 = 1
    50 = 2
    52 = 2
    53 = 2
    55 = 2
    56 = 1
    57 = 1
    58 = 1
    59 = 2
    62 = 1
    63 = 2
    65 = 2
    66 = 1
    70 = 1
    74 = 1
}

foreach ($row in $winners.Keys) {
    $ws.Cells.Item($row, 11).Value = $winners[$row]
}
